# Update the cached "datetimeFigureOut" footer field text from 4/26/2022
# to 4/29/2022 across the slide master and every slide layout.
#
# PowerPoint caches the rendered value of an automatically-updating
# date field inside the <a:fld type="datetimeFigureOut"> run's <a:t>.
# That cached text lives on the master and on each custom layout (the
# slides themselves inherit it), so every one of those placeholders
# needs to be touched.

$p = $ppt.ActivePresentation

$oldDate = "4/26/2022"
$newDate = "4/29/2022"

function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1) Slide master
Update-DatePlaceholder($p.SlideMaster)

# 2) Every slide layout attached to the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder($layouts.Item($li))
}
